$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.281.67"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.92%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.973.33"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.64%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.622"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.79%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.22"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -11.19%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.371"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -6.98%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "56.27"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.88%  "
$ws.Range("E11").Value = "  +7.24%  "
$ws.Range("E12").Value = "  -0.54%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.59"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.32%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.852"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -8.06%  "
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.261.04"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.81%  "
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.83"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -6.68%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.41"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -5.05%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.982.47"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.11%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "36.212.69"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.85%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0879"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.68%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.24"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.41%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "233.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.74%  "
$ws.Range("E24").Value = "  +0.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.50"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.30"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.58%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.82"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.76%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "163.72"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.16%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.133"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.55%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.70"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.67%  "
$ws.Range("E31").Value = "  -2.88%  "
$ws.Range("E32").Value = "  -1.55%  "
$ws.Range("E33").Value = "  -6.59%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0659"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.67%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.45"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.36%  "
$ws.Range("B36").Value = "BinanceUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("B37").Value = "THORChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.11"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.77%  "
$ws.Range("E38").Value = "  -1.48%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.21"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -7.85%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.94"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.47%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.21"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.91%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0960"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.72%  "
$ws.Range("E43").Value = "  -6.31%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0212"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.08%  "
$ws.Range("E45").Value = "  -6.99%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -9.79%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "91.03"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.61%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.371.60"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.52%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.37"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.48%  "
$ws.Range("E50").Value = "  -3.88%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "44.88"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.06%  "
